# Decision_tree with embedded features
#
# Inserts a new row for "Decision Tree Regression+feature selection"
# just above the existing "GB Regression" row, pushing the GB Regression
# row (and its thick-bottom border formatting) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 12; existing row 12 (GB Regression, Id 9) becomes row 13.
$ws.Rows.Item(12).Insert()

# Row 12 currently has the default (no) formatting after the insert -
# clone the formatting from row 11 (same plain side-border style used by
# every interior row of the table) before filling in its values.
$ws.Range("A11:C11").Copy()
$ws.Range("A12:C12").PasteSpecial(-4122)

# New data row: Id 9, "Decision Tree Regression+feature selection".
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Decision Tree Regression+feature selection"
$ws.Range("C12").Value = 85.461516214202803

# The old last row (GB Regression) is now row 13 and becomes Id 10.
$ws.Range("A13").Value = 10

# Reflect the final selection left behind in the saved file.
$ws.Range("C16").Select() | Out-Null
